$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append at the bottom of the existing table (aggiornamento fino a 6/03)
$rows = @(
    @{ r = 245; a = 44319; b = 1; c = 2; d = 62.51953735542357 },
    @{ r = 246; a = 44320; b = 0; c = 2; d = 62.51953735542357 },
    @{ r = 247; a = 44321; b = 0; c = 2; d = 62.51953735542357 }
)

$lastExistingRow = 244

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.a
    $ws.Cells.Item($row.r, 2).Value = $row.b
    $ws.Cells.Item($row.r, 3).Value = $row.c
    $ws.Cells.Item($row.r, 4).Value = $row.d

    # Copy the date cell's formatting (style index 2: centered, bordered, date number format)
    # from the last pre-existing row so the new rows match the rest of the table.
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $ws.Cells.Item($row.r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
